# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker/period table (B16:G21) is re-sorted so that the two workers'
# entries interleave, ordered by "Periodo Mora" ascending (2406, 2407, 2408)
# instead of being grouped by worker with period descending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Worker reference data
$doc1 = "71264684"
$name1 = "ALEJANDRO BERMUDEZ FERNANDEZ"
$doc2 = "45761241"
$name2 = "CLARA INES GUZMAN MARTINEZ"

# New row order (Tipo Doc always "CC"):
#  16: doc1 / name1 / 2406
#  17: doc2 / name2 / 2406
#  18: doc1 / name1 / 2407
#  19: doc2 / name2 / 2407
#  20: doc1 / name1 / 2408
#  21: doc2 / name2 / 2408

$ws.Range("C16").Value = $doc1
$ws.Range("D16").Value = $name1
$ws.Range("E16").Value = "2406"

$ws.Range("C17").Value = $doc2
$ws.Range("D17").Value = $name2
$ws.Range("E17").Value = "2406"

$ws.Range("C18").Value = $doc1
$ws.Range("D18").Value = $name1
$ws.Range("E18").Value = "2407"

$ws.Range("C19").Value = $doc2
$ws.Range("D19").Value = $name2
$ws.Range("E19").Value = "2407"

$ws.Range("C20").Value = $doc1
$ws.Range("D20").Value = $name1
$ws.Range("E20").Value = "2408"

$ws.Range("C21").Value = $doc2
$ws.Range("D21").Value = $name2
$ws.Range("E21").Value = "2408"
